$d = $word.ActiveDocument

$old0 = "M2DocEvaluator.caseConditional(M2DocEvaluator.java:1318)"
$new0 = "M2DocEvaluator.caseConditional(M2DocEvaluator.java:1438)"
$result0 = $d.Content.Find.Execute($old0, $true, $false, $false, $false, $false, $true, 1, $false, $new0, 2)
if (-not $result0) { Write-Host "FAILED replacement 0: $old0" }

$old1 = "M2DocEvaluator.doSwitch(M2DocEvaluator.java:1096)"
$new1 = "M2DocEvaluator.doSwitch(M2DocEvaluator.java:1216)"
$result1 = $d.Content.Find.Execute($old1, $true, $false, $false, $false, $false, $true, 1, $false, $new1, 2)
if (-not $result1) { Write-Host "FAILED replacement 1: $old1" }

$old2 = "M2DocEvaluator.caseBlock(M2DocEvaluator.java:1305)"
$new2 = "M2DocEvaluator.caseBlock(M2DocEvaluator.java:1425)"
$result2 = $d.Content.Find.Execute($old2, $true, $false, $false, $false, $false, $true, 1, $false, $new2, 2)
if (-not $result2) { Write-Host "FAILED replacement 2: $old2" }

$old3 = "M2DocEvaluator.caseConditional(M2DocEvaluator.java:1329)"
$new3 = "M2DocEvaluator.caseConditional(M2DocEvaluator.java:1449)"
$result3 = $d.Content.Find.Execute($old3, $true, $false, $false, $false, $false, $true, 1, $false, $new3, 2)
if (-not $result3) { Write-Host "FAILED replacement 3: $old3" }

$old4 = "M2DocEvaluator.caseDocumentTemplate(M2DocEvaluator.java:283)"
$new4 = "M2DocEvaluator.caseDocumentTemplate(M2DocEvaluator.java:287)"
$result4 = $d.Content.Find.Execute($old4, $true, $false, $false, $false, $false, $true, 1, $false, $new4, 2)
if (-not $result4) { Write-Host "FAILED replacement 4: $old4" }

$old5 = "M2DocEvaluator.generate(M2DocEvaluator.java:272)"
$new5 = "M2DocEvaluator.generate(M2DocEvaluator.java:276)"
$result5 = $d.Content.Find.Execute($old5, $true, $false, $false, $false, $false, $true, 1, $false, $new5, 2)
if (-not $result5) { Write-Host "FAILED replacement 5: $old5" }

$old6 = "AbstractTemplatesTestSuite.prepareoutputAndGenerate(AbstractTemplatesTestSuite.java:479)"
$new6 = "AbstractTemplatesTestSuite.prepareoutputAndGenerate(AbstractTemplatesTestSuite.java:480)"
$result6 = $d.Content.Find.Execute($old6, $true, $false, $false, $false, $false, $true, 1, $false, $new6, 2)
if (-not $result6) { Write-Host "FAILED replacement 6: $old6" }

$old7 = "AbstractTemplatesTestSuite.generation(AbstractTemplatesTestSuite.java:388)"
$new7 = "AbstractTemplatesTestSuite.generation(AbstractTemplatesTestSuite.java:389)"
$result7 = $d.Content.Find.Execute($old7, $true, $false, $false, $false, $false, $true, 1, $false, $new7, 2)
if (-not $result7) { Write-Host "FAILED replacement 7: $old7" }

$old8 = "sun.reflect.GeneratedMethodAccessor75.invoke(Unknown Source)"
$new8 = "sun.reflect.GeneratedMethodAccessor74.invoke(Unknown Source)"
$result8 = $d.Content.Find.Execute($old8, $true, $false, $false, $false, $false, $true, 1, $false, $new8, 2)
if (-not $result8) { Write-Host "FAILED replacement 8: $old8" }

$old9 = "`tat org.eclipse.jdt.internal.junit4.runner.JUnit4TestReference.run(JUnit4TestReference.java:86)`n`tat org.eclipse.jdt.internal.junit.runner.TestExecution.run(TestExecution.java:38)`n`tat org.eclipse.jdt.internal.junit.runner.RemoteTestRunner.runTests(RemoteTestRunner.java:539)`n`tat org.eclipse.jdt.internal.junit.runner.RemoteTestRunner.runTests(RemoteTestRunner.java:761)`n`tat org.eclipse.jdt.internal.junit.runner.RemoteTestRunner.run(RemoteTestRunner.java:461)`n`tat org.eclipse.jdt.internal.junit.runner.RemoteTestRunner.main(RemoteTestRunner.java:207)"
$new9 = "`tat org.apache.maven.surefire.junit4.JUnit4Provider.execute(JUnit4Provider.java:264)`n`tat org.apache.maven.surefire.junit4.JUnit4Provider.executeTestSet(JUnit4Provider.java:153)`n`tat org.apache.maven.surefire.junit4.JUnit4Provider.invoke(JUnit4Provider.java:124)`n`tat sun.reflect.NativeMethodAccessorImpl.invoke0(Native Method)`n`tat sun.reflect.NativeMethodAccessorImpl.invoke(NativeMethodAccessorImpl.java:62)`n`tat sun.reflect.DelegatingMethodAccessorImpl.invoke(DelegatingMethodAccessorImpl.java:43)`n`tat java.lang.reflect.Method.invoke(Method.java:498)`n`tat org.apache.maven.surefire.util.ReflectionUtils.invokeMethodWithArray2(ReflectionUtils.java:208)`n`tat org.apache.maven.surefire.booter.ProviderFactory`$ProviderProxy.invoke(ProviderFactory.java:156)`n`tat org.apache.maven.surefire.booter.ProviderFactory.invokeProvider(ProviderFactory.java:82)`n`tat org.eclipse.tycho.surefire.osgibooter.OsgiSurefireBooter.run(OsgiSurefireBooter.java:91)`n`tat org.eclipse.tycho.surefire.osgibooter.HeadlessTestApplication.run(HeadlessTestApplication.java:21)`n`tat sun.reflect.NativeMethodAccessorImpl.invoke0(Native Method)`n`tat sun.reflect.NativeMethodAccessorImpl.invoke(NativeMethodAccessorImpl.java:62)`n`tat sun.reflect.DelegatingMethodAccessorImpl.invoke(DelegatingMethodAccessorImpl.java:43)`n`tat java.lang.reflect.Method.invoke(Method.java:498)`n`tat org.eclipse.equinox.internal.app.EclipseAppContainer.callMethodWithException(EclipseAppContainer.java:587)`n`tat org.eclipse.equinox.internal.app.EclipseAppHandle.run(EclipseAppHandle.java:198)`n`tat org.eclipse.core.runtime.internal.adaptor.EclipseAppLauncher.runApplication(EclipseAppLauncher.java:134)`n`tat org.eclipse.core.runtime.internal.adaptor.EclipseAppLauncher.start(EclipseAppLauncher.java:104)`n`tat org.eclipse.core.runtime.adaptor.EclipseStarter.run(EclipseStarter.java:388)`n`tat org.eclipse.core.runtime.adaptor.EclipseStarter.run(EclipseStarter.java:243)`n`tat sun.reflect.NativeMethodAccessorImpl.invoke0(Native Method)`n`tat sun.reflect.NativeMethodAccessorImpl.invoke(NativeMethodAccessorImpl.java:62)`n`tat sun.reflect.DelegatingMethodAccessorImpl.invoke(DelegatingMethodAccessorImpl.java:43)`n`tat java.lang.reflect.Method.invoke(Method.java:498)`n`tat org.eclipse.equinox.launcher.Main.invokeFramework(Main.java:656)`n`tat org.eclipse.equinox.launcher.Main.basicRun(Main.java:592)`n`tat org.eclipse.equinox.launcher.Main.run(Main.java:1498)`n`tat org.eclipse.equinox.launcher.Main.main(Main.java:1471)"
$result9 = $d.Content.Find.Execute($old9, $true, $false, $false, $false, $false, $true, 1, $false, $new9, 2)
if (-not $result9) { Write-Host "FAILED replacement 9: $old9" }
